# Commit: "Fruta / hortaliza, semanal"
# The weekly refresh inserts a new daily price record for Jengibre
# (Vega Modelo de Temuco) at the top of the data block (row 43),
# pushing all the following records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 43; existing rows 43:128 shift down to 44:129
$ws.Rows("43:43").Insert()

# Populate the newly inserted row with the latest observation
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 44581
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = 100114007
$ws.Range("G43").Value = "Jengibre"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 40
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("M43").Value = 20000
$ws.Range("N43").Value = "$/caja 13 kilos"
$ws.Range("O43").Value = "Perú"
$ws.Range("P43").Value = 1538
$ws.Range("Q43").Value = 13
$ws.Range("R43").Value = "Hortaliza"
